# Applies the "case with 380 kV done" data update to Sheet1
# (loading_percent.xlsx) - updates computed loading percentages for
# rows 2-25 across columns C, D, E, F, G, H, J, K, O.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2 = @{ "C" = 4.918206080611712; "D" = 10.10000125667195; "E" = 14.05119685983825; "F" = 28.56534670690753; "G" = 27.19863681651143; "H" = 13.98782412387762; "J" = 9.76431838930154; "K" = 16.5798889347404; "O" = 21.06708393494464 }
    3 = @{ "C" = 4.749792681091111; "D" = 10.01835647363327; "E" = 13.97637652016467; "F" = 28.69064401669711; "G" = 27.40496711370968; "H" = 14.06649587242386; "J" = 9.772510979567503; "K" = 15.83430605596674; "O" = 21.20927709565592 }
    4 = @{ "C" = 4.644761917343654; "D" = 9.969201810761561; "E" = 13.932979134866; "F" = 28.777120494908; "G" = 27.54493088166578; "H" = 14.11790600191266; "J" = 9.779172719477804; "K" = 15.35728215786864; "O" = 21.3030494417906 }
    5 = @{ "C" = 4.601629703093863; "D" = 9.94943274609874; "E" = 13.91594767702743; "F" = 28.81474732562095; "G" = 27.60527512432552; "H" = 14.13963647336485; "J" = 9.782297789365748; "K" = 15.1582678979597; "O" = 21.34288268157887 }
    6 = @{ "C" = 4.594449992287688; "D" = 9.946166398445785; "E" = 13.91315946600572; "F" = 28.82113901692212; "G" = 27.61549413316963; "H" = 14.14329192538368; "J" = 9.78284149060319; "K" = 15.12494937081996; "O" = 21.34959466194451 }
    7 = @{ "C" = 4.644181456384703; "D" = 9.96893411644611; "E" = 13.93274677958813; "F" = 28.77761829603401; "G" = 27.54573135224117; "H" = 14.11819590730281; "J" = 9.779213203674399; "K" = 15.35461659560978; "O" = 21.30358009365306 }
    8 = @{ "C" = 4.860518623118351; "D" = 10.0716577104293; "E" = 14.02487974353476; "F" = 28.60656090043778; "G" = 27.26700572493846; "H" = 14.01430535345062; "J" = 9.76680468332157; "K" = 16.32692192864238; "O" = 21.11476697246995 }
    9 = @{ "C" = 5.268793916829893; "D" = 10.28009445003325; "E" = 14.2250721216376; "F" = 28.34739482124816; "G" = 26.82712236608429; "H" = 13.83523597727969; "J" = 9.755408225256247; "K" = 18.07317071913038; "O" = 20.7960515897723 }
    10 = @{ "C" = 5.555352065088186; "D" = 10.43648082791326; "E" = 14.38312226854789; "F" = 28.20421483443814; "G" = 26.57082453208781; "H" = 13.71873507631594; "J" = 9.754905741701736; "K" = 19.24970859808245; "O" = 20.59366142533357 }
    11 = @{ "C" = 5.682152978249169; "D" = 10.50810990494266; "E" = 14.45719184456117; "F" = 28.14947952644094; "G" = 26.46913889092851; "H" = 13.66901336211441; "J" = 9.75638041807265; "K" = 19.7604695716204; "O" = 20.50856456827674 }
    12 = @{ "C" = 5.729611496297525; "D" = 10.53528550327614; "E" = 14.48553346063061; "F" = 28.13025898917184; "G" = 26.43280788334946; "H" = 13.65065679923849; "J" = 9.757183045670992; "K" = 19.95027588934249; "O" = 20.47734982440641 }
    13 = @{ "C" = 5.719415972042128; "D" = 10.52943077534561; "E" = 14.47941686312838; "F" = 28.13433129887238; "G" = 26.44053517714891; "H" = 13.65458920570329; "J" = 9.756999337538007; "K" = 19.90955945173291; "O" = 20.4840274615885 }
    14 = @{ "C" = 5.686068856270103; "D" = 10.51034473969863; "E" = 14.45951776454029; "F" = 28.14786798552398; "G" = 26.46610614442573; "H" = 13.66749369162794; "J" = 9.75644156121824; "K" = 19.77615782564579; "O" = 20.50597623763793 }
    15 = @{ "C" = 5.665568787354171; "D" = 10.4986601000972; "E" = 14.44736658025571; "F" = 28.15635611728463; "G" = 26.4820532990981; "H" = 13.67545955748272; "J" = 9.756131685381057; "K" = 19.69397318984755; "O" = 20.51955218833463 }
    16 = @{ "C" = 5.54698963862323; "D" = 10.43180789916391; "E" = 14.37832371740534; "F" = 28.20800213174221; "G" = 26.57777259194616; "H" = 13.72205049636822; "J" = 9.75484358281685; "K" = 19.21582893170222; "O" = 20.59936359228321 }
    17 = @{ "C" = 5.473299295605638; "D" = 10.39090775457144; "E" = 14.33651088605362; "F" = 28.24235671724629; "G" = 26.6403337073851; "H" = 13.75147206855612; "J" = 9.754489142135879; "K" = 18.91616980236812; "O" = 20.65011586586546 }
    18 = @{ "C" = 5.430582308569727; "D" = 10.36743015378333; "E" = 14.31266685410074; "F" = 28.26309463596109; "G" = 26.67771831980563; "H" = 13.76870278518806; "J" = 9.754445647509225; "K" = 18.74151662620957; "O" = 20.67996298243454 }
    19 = @{ "C" = 5.416063526559973; "D" = 10.35948970624016; "E" = 14.30462956442151; "F" = 28.27028376808606; "G" = 26.69061577326134; "H" = 13.7745897147137; "J" = 9.754458484796711; "K" = 18.68199057746601; "O" = 20.69018111513174 }
    20 = @{ "C" = 5.481178502429023; "D" = 10.39525690967212; "E" = 14.34094079188932; "F" = 28.23859830236089; "G" = 26.63352872800808; "H" = 13.74830818659803; "J" = 9.754510280091095; "K" = 18.94830749099536; "O" = 20.64464528088693 }
    21 = @{ "C" = 5.695879209678833; "D" = 10.51594953181912; "E" = 14.46535481378471; "F" = 28.14385095858036; "G" = 26.4585360489299; "H" = 13.66369051378196; "J" = 9.756598772558226; "K" = 19.81543968302999; "O" = 20.49950189277241 }
    22 = @{ "C" = 5.832926719247244; "D" = 10.59511894832453; "E" = 14.54836559563935; "F" = 28.09071412069977; "G" = 26.3568608550505; "H" = 13.61113949423654; "J" = 9.759386678595948; "K" = 20.36109997869378; "O" = 20.4105303544494 }
    23 = @{ "C" = 5.760095360529668; "D" = 10.55284447330691; "E" = 14.50391218968238; "F" = 28.11826675929236; "G" = 26.40995511910209; "H" = 13.63893484831177; "J" = 9.757768793148294; "K" = 20.07182415947552; "O" = 20.45747499681546 }
    24 = @{ "C" = 5.477617404048452; "D" = 10.39329054172383; "E" = 14.33893742348516; "F" = 28.24029440808193; "G" = 26.63660084644783; "H" = 13.74973759378157; "J" = 9.754500224316843; "K" = 18.9337854329456; "O" = 20.64711644907507 }
    25 = @{ "C" = 5.160456844299305; "D" = 10.22306187677979; "E" = 14.16891812722361; "F" = 28.40926292874204; "G" = 26.93450890247578; "H" = 13.88103629899967; "J" = 9.757107497407041; "K" = 17.61892371273724; "O" = 20.87671896000701 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
